# Adds the results of a new weekly "pelada" (game session) as rows 130-152
# of the Jogadores sheet, appending to the existing running totals table.
# Columns: A=Jogadores, C=Vitorias, D=Empate, E=Derrotas, F=Gols,
#          G=Partidas, H=Tarde de Vitoria, I=La barca, J=Craque do Dia,
#          K=Gols Sofridos  (column B=Pontos is intentionally left blank,
#          matching every prior block in this sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
  @(130, 'Athos', 2,1,3,2,1,0,1,0,0),
  @(131, 'Miqueias', 2,1,3,1,1,0,1,0,0),
  @(132, 'Boneco', 2,1,3,0,1,0,1,0,0),
  @(133, 'Marcos', 2,1,3,0,1,0,1,0,0),
  @(134, 'Corinthiano', 2,1,3,0,1,0,1,0,0),
  @(135, 'Marcelão', 2,2,2,1,1,0,0,0,0),
  @(136, 'Jorge', 2,2,2,1,1,0,0,0,0),
  @(137, 'Fernando', 2,2,2,2,1,0,0,1,0),
  @(138, 'Digão', 2,2,2,0,1,0,0,0,0),
  @(139, 'Juscielio', 2,2,2,0,1,0,0,0,0),
  @(140, 'Romario', 2,2,1,1,1,1,0,0,0),
  @(141, 'Eduardo', 2,2,1,2,1,1,0,0,0),
  @(142, 'Leandrão', 2,2,1,0,1,1,0,0,0),
  @(143, 'Cabeleira', 2,2,1,3,1,1,0,0,0),
  @(144, 'Leandrinho', 2,2,1,0,1,1,0,0,0),
  @(145, 'Peixe', 1,3,1,0,1,0,0,0,0),
  @(146, 'Heider', 1,3,1,0,1,0,0,0,0),
  @(147, 'Ismael', 1,3,1,1,1,0,0,0,0),
  @(148, 'David', 1,3,1,1,1,0,0,0,0),
  @(149, 'Eder', 1,3,1,1,1,0,0,0,0),
  @(150, 'Matheus', 3,1,3,0,1,0,0,0,8),
  @(151, 'Igor Goleiro', 1,4,2,0,1,0,1,0,5),
  @(152, 'Chelin', 3,3,1,0,1,1,0,0,3)
)

foreach ($row in $newRows) {
  $r = $row[0]
  $ws.Cells.Item($r, 1).Value = $row[1]
  $ws.Range("C$r").Value = $row[2]
  $ws.Range("D$r").Value = $row[3]
  $ws.Range("E$r").Value = $row[4]
  $ws.Range("F$r").Value = $row[5]
  $ws.Range("G$r").Value = $row[6]
  $ws.Range("H$r").Value = $row[7]
  $ws.Range("I$r").Value = $row[8]
  $ws.Range("J$r").Value = $row[9]
  $ws.Range("K$r").Value = $row[10]
}

# Move the frozen-pane scroll position and active selection to mirror
# where the author ended up after appending the new block.
$win = $excel.ActiveWindow
$win.ScrollRow = 145
$win.ScrollColumn = 1
$ws.Range("C153").Select()
